$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update H6: was a plain value (1000000), becomes a formula 1000000+20000 (=1020000)
$ws.Range("H6").Formula = "=1000000+20000"

# Update H9: was a plain value (1620000), becomes a formula 1000000+20000+620000 (=1640000)
$ws.Range("H9").Formula = "=1000000+20000+620000"

# Update the selection shown in the sheet view (activeCell H9 -> H6)
$ws.Range("H6").Select()

# Update the workbook window position (xWindow/yWindow) to match new values
$excel.ActiveWindow.Left = 790
$excel.ActiveWindow.Top = 960

$wb.Save()
